$d = $word.ActiveDocument

# 1) Merge the split "i" run (inside proofErr spell-check markers) back into
#    plain flowing text by re-finding/replacing the surrounding phrase with
#    itself. Word's Find/Replace rebuilds the run(s) it touches as a single
#    run and drops the now-redundant proofErr bookmarks.
$needle = "doesn" + [char]0x2019 + "t can i still call"
$d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)

# 2) Append the new development-log paragraphs after "Need to tidy the
#    project", before the trailing blank paragraph, using InsertXML on a
#    collapsed range so no existing content is disturbed.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Need to tidy the project") {
        $target = $cand
    }
}

$insertAt = $target.Range.End
$insertRange = $d.Range($insertAt, $insertAt)

$apos = [char]0x2019
$newBodyXml = "<w:p/>" + `
    "<w:p><w:r><w:t>Look into view port rectangle forsplit screen</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>Take triggers off colliders to make the boat collide</w:t></w:r>" + `
    "<w:r><w:t xml:space=`"preserve`"> freez x z rotation and y position but seems not to turn or very little, maybe because so little difference in angle force being applied need to check though think it is 30 degree</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>Colliders strangley push up triggers too on the camera triggers</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>Need change code so doesn${apos}t accelerate out of hand</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>But then issues of spinning and sinking and then not like turning once freezing directions</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>Also the speed now accelerates probably because when not a trigger it suddenly has weight</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>And so now it builds momentum</w:t></w:r></w:p>" + `
    "<w:p><w:r><w:t>Needs more speed to move too</w:t></w:r></w:p>" + `
    "<w:p/>"

$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    "<w:body>$newBodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$insertRange.InsertXML($pkgXml)
